# Remove incubation-related variables (rows) from the "vocabulary" sheet:
#   row 12: incubation_protocol
#   row 13: nylon_pore_size_µm
#   row 14: bag_size_cm2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vocabulary")

# Delete the three rows entirely (bottom-up so row numbers stay valid).
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()

# Update the frozen-pane / selection state to match the post-edit view.
$ws.Activate()
$ws.Range("M2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("M20").Select()
